# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" wherever it appears
# - Shrink the "Status"-related columns to match the regenerated report's autosized widths

$wb = $excel.ActiveWorkbook

# --- 1. Replace the status text across every sheet/cell that currently holds it ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $cellText = [string]$cell.Value2
            if ($cellText -eq "Ready for handoff") {
                $cell.Value2 = "In Translation"
            }
        }
    }
}

# --- 2. Narrow the status columns' width (regenerated report reflows them) ---
# Overview sheet: columns E ("zh-cn") and F ("de-de") hold the status values
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C holds "Status"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C holds "Status"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
